# Mise a jour bd
# Adds two new logbook entries ("12:40" and "16 :15") after the existing
# "11:40 : Je viens de terminer ..." entry in the "Déroulement" section of
# J2, and moves the trailing _GoBack bookmark to the end of the new last
# paragraph (mirroring what Word does when new text is typed at the
# previous cursor/_GoBack location).

$d = $word.ActiveDocument

$nbsp = [char]0x00A0
$cr = [char]13

# The "_GoBack" bookmark currently sits right at the end of the "11:40"
# paragraph (after "... après." and before the paragraph mark). Remove it
# first so our insertion below lands cleanly at the end of that paragraph's
# text instead of being pushed after/around the bookmark.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the end of the "11:40 : ..." entry via its closing sentence.
$target = $d.Content
$found = $target.Find.Execute("de données après.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor paragraph text"
}
$target.Collapse(0)

$entry1 = "12:40${nbsp}: Je me remets a travailler directement. Je vais continuer de documenter la base de données ainsi que toutes les classes containers que j’utiliserai plus tard dans le code."
$entry2 = "16 :15${nbsp}: J’ai bien avancé dans l’explication des classes conteneurs et j’ai terminer l’explication des différentes tables de la base de données. Je mets en pause mon travail pour pouvoir rejoindre le rendez-vous journalier, sur Meet, avec mon maître d’apprentissage."

$block = "$cr$cr$entry1$cr$cr$entry2"
$target.InsertAfter($block)

# $target now ends exactly at the paragraph mark that follows the new last
# paragraph ("16 :15 : ... apprentissage."). Re-create the _GoBack bookmark
# there. A zero-length bookmark added exactly at a paragraph-mark position
# is mishandled by this COM host, so: temporarily insert a placeholder
# character past that point, add the bookmark while it is safely before
# that placeholder (not a paragraph boundary), then remove the
# placeholder again — the bookmark stays correctly anchored in place.
$endPos = $target.End
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$junk = $d.Range($endPos, $endPos + 1)
$junk.Delete()
